$d = $word.ActiveDocument

$d.Content.Find.Execute("91÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷2=", 2) | Out-Null
$d.Content.Find.Execute("95÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷4=", 2) | Out-Null
$d.Content.Find.Execute("70÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷7=", 2) | Out-Null
$d.Content.Find.Execute("56÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷8=", 2) | Out-Null
$d.Content.Find.Execute("41÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷6=", 2) | Out-Null
$d.Content.Find.Execute("50÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷2=", 2) | Out-Null
$d.Content.Find.Execute("99÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷9=", 2) | Out-Null
$d.Content.Find.Execute("53÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷3=", 2) | Out-Null
$d.Content.Find.Execute("90÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷4=", 2) | Out-Null
$d.Content.Find.Execute("50÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷6=", 2) | Out-Null
$d.Content.Find.Execute("74÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷7=", 2) | Out-Null
$d.Content.Find.Execute("83÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷5=", 2) | Out-Null
$d.Content.Find.Execute("88÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷4=", 2) | Out-Null
$d.Content.Find.Execute("11÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷9=", 2) | Out-Null
$d.Content.Find.Execute("47÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷2=", 2) | Out-Null
$d.Content.Find.Execute("68÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷7=", 2) | Out-Null
$d.Content.Find.Execute("79÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷9=", 2) | Out-Null
$d.Content.Find.Execute("11÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷5=", 2) | Out-Null
$d.Content.Find.Execute("72÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷7=", 2) | Out-Null
$d.Content.Find.Execute("72÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=", 2) | Out-Null
$d.Content.Find.Execute("40÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷4=", 2) | Out-Null
$d.Content.Find.Execute("71÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷4=", 2) | Out-Null
$d.Content.Find.Execute("24÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷2=", 2) | Out-Null
$d.Content.Find.Execute("20÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷3=", 2) | Out-Null
$d.Content.Find.Execute("89÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷5=", 2) | Out-Null
